# Weekly fruit/vegetable price update: a new weekly price record for
# "Berenjena" (Vega Modelo de Temuco) is inserted at row 171, pushing the
# existing rows 171:198 down to 172:199 (dimension grows from R198 to R199).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 171; this shifts rows 171-198
# down to 172-199 and extends the sheet dimension automatically.
$ws.Rows(171).Insert()

# Populate the newly inserted row 171 with this week's record.
$ws.Cells.Item(171, 1).Value  = 10
$ws.Cells.Item(171, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(171, 3).Value  = "La Araucanía"
$ws.Cells.Item(171, 4).Value  = 44505
$ws.Cells.Item(171, 5).Value  = 9
$ws.Cells.Item(171, 6).Value  = 100112001
$ws.Cells.Item(171, 7).Value  = "Berenjena"
$ws.Cells.Item(171, 8).Value  = "Sin especificar"
$ws.Cells.Item(171, 9).Value  = "Primera"
$ws.Cells.Item(171, 10).Value = 55
$ws.Cells.Item(171, 11).Value = 12000
$ws.Cells.Item(171, 12).Value = 12000
$ws.Cells.Item(171, 13).Value = 12000
$ws.Cells.Item(171, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(171, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(171, 16).Value = 200
$ws.Cells.Item(171, 17).Value = 60
$ws.Cells.Item(171, 18).Value = "Hortaliza"
